$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 24260
$ws.Range("E2").Value = 2715
$ws.Range("F2").Value = 2715
$ws.Range("G2").Value = 1873
$ws.Range("H2").Value = 1299
$ws.Range("I2").Value = 1273
$ws.Range("J2").Value = 26
$ws.Range("K2").Value = 53070
$ws.Range("L2").Value = 27047
$ws.Range("M2").Value = 26024
$ws.Range("N2").Value = 25854
$ws.Range("O2").Value = 169
$ws.Range("P2").Value = 100
$ws.Range("Q2").Value = 1850
$ws.Range("R2").Value = -1007
$ws.Range("S2").Value = -975
$ws.Range("T2").Value = 1540
$ws.Range("U2").Value = 310
$ws.Range("V2").Value = 15899
$ws.Range("W2").Value = 11.19
$ws.Range("X2").Value = 5.36
$ws.Range("Y2").Value = 5.01
$ws.Range("Z2").Value = 2.45
$ws.Range("AA2").Value = 103.93
$ws.Range("AB2").Value = 32598.13
$ws.Range("AC2").Value = 2546
$ws.Range("AD2").Value = 83.86
$ws.Range("AE2").Value = 58760
$ws.Range("AF2").Value = 3.63
$ws.Range("AG2").Value = 2000
$ws.Range("AH2").Value = 0.9399999999999999
$ws.Range("AI2").Value = 69.13
$ws.Range("AJ2").Value = 50000000

# Row 3
$ws.Range("D3").Value = 392995
$ws.Range("E3").Value = 14068
$ws.Range("F3").Value = 14033
$ws.Range("G3").Value = 62881
$ws.Range("H3").Value = 55493
$ws.Range("I3").Value = 53460
$ws.Range("J3").Value = 2033
$ws.Range("K3").Value = 966338
$ws.Range("L3").Value = 565566
$ws.Range("M3").Value = 400772
$ws.Range("N3").Value = 125665
$ws.Range("O3").Value = 275107
$ws.Range("P3").Value = 154
$ws.Range("Q3").Value = 45889
$ws.Range("R3").Value = 43299
$ws.Range("S3").Value = -22255
$ws.Range("T3").Value = 22377
$ws.Range("U3").Value = 23512
$ws.Range("V3").Value = 317217
$ws.Range("W3").Value = 3.58
$ws.Range("X3").Value = 14.12
$ws.Range("Y3").Value = 70.56
$ws.Range("Z3").Value = 10.89
$ws.Range("AA3").Value = 141.12
$ws.Range("AB3").Value = 87618.89999999999
$ws.Range("AC3").Value = 93713
$ws.Range("AD3").Value = 2.57
$ws.Range("AE3").Value = 222855
$ws.Range("AF3").Value = 1.08
$ws.Range("AG3").Value = 3400
$ws.Range("AH3").Value = 1.41
$ws.Range("AI3").Value = 3.59
$ws.Range("AJ3").Value = 70360297

# Row 4
$ws.Range("D4").Value = 827299
$ws.Range("E4").Value = 52808
$ws.Range("F4").Value = 52982
$ws.Range("G4").Value = 43897
$ws.Range("H4").Value = 28225
$ws.Range("I4").Value = 7656
$ws.Range("J4").Value = 20569
$ws.Range("K4").Value = 1030473
$ws.Range("L4").Value = 607206
$ws.Range("M4").Value = 423267
$ws.Range("N4").Value = 131124
$ws.Range("O4").Value = 292144
$ws.Range("P4").Value = 154
$ws.Range("Q4").Value = 91107
$ws.Range("R4").Value = -73805
$ws.Range("S4").Value = -16224
$ws.Range("T4").Value = 49354
$ws.Range("U4").Value = 41753
$ws.Range("V4").Value = 313808
$ws.Range("W4").Value = 6.38
$ws.Range("X4").Value = 3.41
$ws.Range("Y4").Value = 5.96
$ws.Range("Z4").Value = 2.83
$ws.Range("AA4").Value = 143.46
$ws.Range("AB4").Value = 91087.10000000001
$ws.Range("AC4").Value = 10794
$ws.Range("AD4").Value = 21.26
$ws.Range("AE4").Value = 232535
$ws.Range("AF4").Value = 0.99
$ws.Range("AG4").Value = 3700
$ws.Range("AH4").Value = 1.61
$ws.Range("AI4").Value = 27.26
$ws.Range("AJ4").Value = 70360297

# Row 5
$ws.Range("D5").Value = 906128
$ws.Range("E5").Value = 57450
$ws.Range("F5").Value = 57450
$ws.Range("G5").Value = 74431
$ws.Range("H5").Value = 50665
$ws.Range("I5").Value = 16774
$ws.Range("J5").Value = 33890
$ws.Range("K5").Value = 1097835
$ws.Range("L5").Value = 640270
$ws.Range("M5").Value = 457564
$ws.Range("N5").Value = 136601
$ws.Range("O5").Value = 320963
$ws.Range("P5").Value = 154
$ws.Range("Q5").Value = 69301
$ws.Range("R5").Value = -62929
$ws.Range("S5").Value = -4080
$ws.Range("T5").Value = 49719
$ws.Range("U5").Value = 19582
$ws.Range("V5").Value = 329084
$ws.Range("W5").Value = 6.34
$ws.Range("X5").Value = 5.59
$ws.Range("Y5").Value = 12.53
$ws.Range("Z5").Value = 4.76
$ws.Range("AA5").Value = 139.93
$ws.Range("AB5").Value = 96552.39
$ws.Range("AC5").Value = 23650
$ws.Range("AD5").Value = 11.97
$ws.Range("AE5").Value = 242249
$ws.Range("AF5").Value = 1.17
$ws.Range("AG5").Value = 4000
$ws.Range("AH5").Value = 1.41
$ws.Range("AI5").Value = 13.45
$ws.Range("AJ5").Value = 70360297

# Row 6
$ws.Range("D6").Value = 1015020
$ws.Range("E6").Value = 46881
$ws.Range("F6").Value = 46881
$ws.Range("G6").Value = 75284
$ws.Range("H6").Value = 61511
$ws.Range("I6").Value = 22531
$ws.Range("K6").Value = 1194570
$ws.Range("L6").Value = 685691
$ws.Range("M6").Value = 508879
$ws.Range("N6").Value = 166275
$ws.Range("P6").Value = 154
$ws.Range("Q6").Value = 78565
$ws.Range("R6").Value = -101818
$ws.Range("S6").Value = 19297
$ws.Range("T6").Value = 62756
$ws.Range("U6").Value = 15809
$ws.Range("V6").Value = 373631
$ws.Range("W6").Value = 4.62
$ws.Range("X6").Value = 6.06
$ws.Range("Y6").Value = 14.88
$ws.Range("Z6").Value = 5.37
$ws.Range("AA6").Value = 134.75
$ws.Range("AB6").Value = 115396.6
$ws.Range("AC6").Value = 31767
$ws.Range("AD6").Value = 8.18
$ws.Range("AE6").Value = 294873
$ws.Range("AF6").Value = 0.88
$ws.Range("AG6").Value = 5000
$ws.Range("AH6").Value = 1.92
$ws.Range("AI6").Value = 12.51
$ws.Range("AJ6").Value = 70360297

# Row 7
$ws.Range("D7").Value = 1004260
$ws.Range("E7").Value = 46033
$ws.Range("G7").Value = 45174
$ws.Range("H7").Value = 32605
$ws.Range("I7").Value = 13849
$ws.Range("K7").Value = 1264095
$ws.Range("L7").Value = 718093
$ws.Range("M7").Value = 546001
$ws.Range("N7").Value = 180281
$ws.Range("P7").Value = 152
$ws.Range("Q7").Value = 85609
$ws.Range("R7").Value = -96439
$ws.Range("S7").Value = 16827
$ws.Range("T7").Value = 69516
$ws.Range("U7").Value = 12526
$ws.Range("W7").Value = 4.58
$ws.Range("X7").Value = 3.25
$ws.Range("Y7").Value = 7.99
$ws.Range("Z7").Value = 2.65
$ws.Range("AA7").Value = 131.52
$ws.Range("AC7").Value = 19526
$ws.Range("AD7").Value = 12.06
$ws.Range("AE7").Value = 340998
$ws.Range("AF7").Value = 0.6899999999999999
$ws.Range("AG7").Value = 5273
$ws.Range("AH7").Value = 2.24
$ws.Range("AI7").Value = 26.79

# Row 8
$ws.Range("D8").Value = 1033904
$ws.Range("E8").Value = 54030
$ws.Range("G8").Value = 56290
$ws.Range("H8").Value = 40851
$ws.Range("I8").Value = 16876
$ws.Range("K8").Value = 1307088
$ws.Range("L8").Value = 724609
$ws.Range("M8").Value = 582480
$ws.Range("N8").Value = 194748
$ws.Range("P8").Value = 152
$ws.Range("Q8").Value = 86822
$ws.Range("R8").Value = -72316
$ws.Range("S8").Value = -13355
$ws.Range("T8").Value = 56262
$ws.Range("U8").Value = 38186
$ws.Range("W8").Value = 5.23
$ws.Range("X8").Value = 3.95
$ws.Range("Y8").Value = 9
$ws.Range("Z8").Value = 3.18
$ws.Range("AA8").Value = 124.4
$ws.Range("AC8").Value = 23794
$ws.Range("AD8").Value = 9.9
$ws.Range("AE8").Value = 368362
$ws.Range("AF8").Value = 0.64
$ws.Range("AG8").Value = 5591
$ws.Range("AH8").Value = 2.37
$ws.Range("AI8").Value = 23.31

# Row 9
$ws.Range("D9").Value = 1068779
$ws.Range("E9").Value = 57458
$ws.Range("G9").Value = 61493
$ws.Range("H9").Value = 44568
$ws.Range("I9").Value = 19167
$ws.Range("K9").Value = 1347032
$ws.Range("L9").Value = 725453
$ws.Range("M9").Value = 621578
$ws.Range("N9").Value = 209942
$ws.Range("P9").Value = 152
$ws.Range("Q9").Value = 86561
$ws.Range("R9").Value = -74974
$ws.Range("S9").Value = -13287
$ws.Range("T9").Value = 56860
$ws.Range("U9").Value = 35141
$ws.Range("W9").Value = 5.38
$ws.Range("X9").Value = 4.17
$ws.Range("Y9").Value = 9.470000000000001
$ws.Range("Z9").Value = 3.36
$ws.Range("AA9").Value = 116.71
$ws.Range("AC9").Value = 27024
$ws.Range("AD9").Value = 8.710000000000001
$ws.Range("AE9").Value = 397101
$ws.Range("AF9").Value = 0.59
$ws.Range("AG9").Value = 6000
$ws.Range("AH9").Value = 2.55
$ws.Range("AI9").Value = 22.02
